$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column D ("q") for rows 5 and 6
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 2

# Update the active selection to D5
$ws.Range("D5").Select()
